$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: phone "09876543" (keep as text so the leading zero survives),
# birthday left blank, total_points reset to 0.
$ws.Range("A12").Value = "'09876543"
$ws.Range("B12").Font.Bold = $false
$ws.Range("C12").Value = 0
